# Edit LOT2051.xlsx: update Objetivos text, insert a row for "Docentes responsaveis" value,
# update Programa resumido / Programa / Metodo / Criterio / Norma de recuperacao / Bibliografia content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Objetivos:" (row 10) English->full objective text; the old faculty name value
#    that used to sit in B10/C10 will be relocated to a new row inserted below row 12.
$ws.Range("B10").Value = 'Fornecer ao aluno os conhecimentos fundamentais relativos ao cultivo de células animais, visando seu emprego como instrumento de obtenção de produtos biotecnológicos de alto valor agregado'
$ws.Range("C10").Value = 'Fornecer ao aluno os conhecimentos fundamentais relativos ao cultivo de células animais, visando seu emprego como instrumento de obtenção de produtos biotecnológicos de alto valor agregado'

# 2. Insert a new row at position 13 (pushes old rows 13-23 down to 14-24), to hold the
#    "Docentes responsaveis:" value (faculty name) in columns B/C, with no A value.
$ws.Rows("13:13").Insert()
$ws.Range("A13").Clear()

# Copy the B/C number format & style from row 10 (same visual style family) into B13:C13
# so the inserted cells pick up style index 2 / 3 like their neighbours.
$ws.Range("B10:C10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B13").Value = '101761 - Arnaldo Márcio Ramalho Prata'
$ws.Range("C13").Value = '101761 - Arnaldo Márcio Ramalho Prata'

# 3. Row 14 (was row 13, "Programa resumido:") -- replace short-summary text (pt-BR)
$ws.Range("B14").Value = 'Introdução à Tecnologia de Cultivo de Células Animais, Crescimento e Morte de Células Animais Cultivadas in vitro, Biorreatores para Células Animais, Aplicações do Cultivo de Células Animais'
$ws.Range("C14").Value = 'Introdução à Tecnologia de Cultivo de Células Animais, Crescimento e Morte de Células Animais Cultivadas in vitro, Biorreatores para Células Animais, Aplicações do Cultivo de Células Animais'

# 4. Row 16 (was row 15, "Programa:") -- replace with the full pt-BR programme text
$ws.Range("B16").Value = '1. Introdução à Tecnologia de Cultivo de Células Animais  Principais marcos e razões da cultura de células animais, Tipos de culturas de células animais, Emprego de células animais. 2. Mecanismo de Crescimento e Morte de Células Animais Cultivadas in vitro  Mecanismos de proliferação celular, Mecanismos de morte celular, Influência das condições ambientais sobre a morte celular, Métodos de detecção da morte celular, Controle da apoptose por técnicas moleculares. 3. Biorreatores para Células Animais  Propagação de inóculo e sistemas de cultivo em pequena escala, Tipos de biorreatores, Aeração e agitação, Aspectos econômicos na seleção de biorreatores. 4. Aplicações do Cultivo de Células Animais  Proteinas recombinantes terapêuticas, Anticorpos monoclonais, Vacinas virais, Bioinseticidas, Terapias celulares e células-tronco, Terapia gênica'
$ws.Range("C16").Value = '1. Introdução à Tecnologia de Cultivo de Células Animais  Principais marcos e razões da cultura de células animais, Tipos de culturas de células animais, Emprego de células animais. 2. Mecanismo de Crescimento e Morte de Células Animais Cultivadas in vitro  Mecanismos de proliferação celular, Mecanismos de morte celular, Influência das condições ambientais sobre a morte celular, Métodos de detecção da morte celular, Controle da apoptose por técnicas moleculares. 3. Biorreatores para Células Animais  Propagação de inóculo e sistemas de cultivo em pequena escala, Tipos de biorreatores, Aeração e agitação, Aspectos econômicos na seleção de biorreatores. 4. Aplicações do Cultivo de Células Animais  Proteinas recombinantes terapêuticas, Anticorpos monoclonais, Vacinas virais, Bioinseticidas, Terapias celulares e células-tronco, Terapia gênica'

# 5. Row 19 (was row 18, "Metodo:") -- replace wrongly-duplicated faculty text with the
#    evaluation method text.
$ws.Range("B19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Range("C19").Value = 'A avaliação será feita por meio de provas escritas.'

# 6. Row 20 (was row 19, "Criterio:") -- replace with the final-grade formula text.
$ws.Range("B20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = [P1 +(2 x P2)] / 3'
$ws.Range("C20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = [P1 +(2 x P2)] / 3'

# 7. Row 21 (was row 20, "Norma de recuperacao:") -- replace with the make-up exam text.
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'

# 8. Row 22 (was row 21, "Bibliografia:") -- replace with the full bibliography text.
$ws.Range("B22").Value = '1. MORAES, A.M., AUGUSTO, E.F.P., CASTILHO, L.R. Tecnologia do Cultivo de Células Animais – de Biofármacos a Terapia Gênica. São Paulo: Rocca, 2008.2. VITOLO, M. (Coordenador). Biotecnologia Farmacêutica – Aspectos sobre aplicação industrial. São Paulo: Edgard Blücher Ltda, 2015.3. SHULER, M.L., KARGI, F. Bioprocess Engineering – Basic Concepts. Second edition. New Jersey: Prentice Hall, 2002.'
$ws.Range("C22").Value = '1. MORAES, A.M., AUGUSTO, E.F.P., CASTILHO, L.R. Tecnologia do Cultivo de Células Animais – de Biofármacos a Terapia Gênica. São Paulo: Rocca, 2008.2. VITOLO, M. (Coordenador). Biotecnologia Farmacêutica – Aspectos sobre aplicação industrial. São Paulo: Edgard Blücher Ltda, 2015.3. SHULER, M.L., KARGI, F. Bioprocess Engineering – Basic Concepts. Second edition. New Jersey: Prentice Hall, 2002.'
